$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal TEXT value even when the text looks numeric
# (e.g. "10"), by briefly marking the cell as Text-formatted, then restoring
# the Normal style so no visible formatting change remains.
function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 19
$ws.Range("A19").Value = 111927215
$ws.Range("B19").Value = 90792
Set-TextValue "I19" '10'
Set-TextValue "J19" 'fruktkroppar'
$ws.Range("Q19").Value = 663486
$ws.Range("R19").Value = 6602647

# Row 20
$ws.Range("B20").Value = 90792

# Row 21
$ws.Range("B21").Value = 90789

# Row 22
$ws.Range("A22").Value = 111926769
$ws.Range("B22").Value = 90792
$ws.Range("I22").ClearContents()
$ws.Range("J22").ClearContents()
$ws.Range("Q22").Value = 663476
$ws.Range("R22").Value = 6602651

# Row 23
$ws.Range("A23").Value = 112083958
$ws.Range("B23").Value = 98961
Set-TextValue "I23" '10'
Set-TextValue "J23" 'plantor/tuvor'
Set-TextValue "K23" 'fullt utvecklade blad'
$ws.Range("Q23").Value = 663551
$ws.Range("R23").Value = 6602700
$ws.Range("S23").Value = 5

# Row 24
$ws.Range("A24").Value = 112083905
$ws.Range("B24").Value = 98961
Set-TextValue "I24" '400'
Set-TextValue "J24" 'stjälkar/strån/skott'
$ws.Range("Q24").Value = 663568
$ws.Range("R24").Value = 6602721

# Row 25
$ws.Range("B25").Value = 98961

# Row 26
$ws.Range("B26").Value = 89033

# Row 27
$ws.Range("A27").Value = 112083737
$ws.Range("B27").Value = 98961
Set-TextValue "I27" '200'
Set-TextValue "J27" 'stjälkar/strån/skott'
$ws.Range("Q27").Value = 663545
$ws.Range("R27").Value = 6602752
$ws.Range("S27").Value = 30

# Row 28
$ws.Range("A28").Value = 112083991
$ws.Range("B28").Value = 98961
Set-TextValue "I28" '300'
$ws.Range("Q28").Value = 663568
$ws.Range("R28").Value = 6602664
$ws.Range("S28").Value = 10

# Row 29
$ws.Range("A29").Value = 112084040
$ws.Range("B29").Value = 98961
$ws.Range("I29").ClearContents()
$ws.Range("J29").ClearContents()
$ws.Range("K29").ClearContents()
$ws.Range("Q29").Value = 663585
$ws.Range("R29").Value = 6602704
$ws.Range("S29").Value = 10

# Row 30
$ws.Range("A30").Value = 112083804
$ws.Range("B30").Value = 98961
$ws.Range("Q30").Value = 663572
$ws.Range("R30").Value = 6602738
